$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(8, 9).Value = 'b'
$ws.Cells.Item(8, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(19, 9).Value = 'sv'
$ws.Cells.Item(19, 10).Value = 'Statement-opinion'
$ws.Cells.Item(23, 9).Value = 'sv'
$ws.Cells.Item(23, 10).Value = 'Statement-opinion'
$ws.Cells.Item(25, 9).Value = 'b'
$ws.Cells.Item(25, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(32, 9).Value = 'sv'
$ws.Cells.Item(32, 10).Value = 'Statement-opinion'
$ws.Cells.Item(33, 9).Value = 'sd'
$ws.Cells.Item(33, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(42, 9).Value = 'sd'
$ws.Cells.Item(42, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(47, 9).Value = 'aa'
$ws.Cells.Item(47, 10).Value = 'Agree/Accept'
$ws.Cells.Item(70, 9).Value = 'sv'
$ws.Cells.Item(70, 10).Value = 'Statement-opinion'
$ws.Cells.Item(71, 9).Value = 'sd'
$ws.Cells.Item(71, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(75, 9).Value = 'sd'
$ws.Cells.Item(75, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(79, 9).Value = 'sd'
$ws.Cells.Item(79, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(80, 9).Value = 'sd'
$ws.Cells.Item(80, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(81, 9).Value = 'aa'
$ws.Cells.Item(81, 10).Value = 'Agree/Accept'
$ws.Cells.Item(87, 9).Value = 'aa'
$ws.Cells.Item(87, 10).Value = 'Agree/Accept'
$ws.Cells.Item(88, 9).Value = 'sv'
$ws.Cells.Item(88, 10).Value = 'Statement-opinion'
$ws.Cells.Item(95, 9).Value = 'aa'
$ws.Cells.Item(95, 10).Value = 'Agree/Accept'
$ws.Cells.Item(96, 9).Value = 'sd'
$ws.Cells.Item(96, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(107, 9).Value = 'sd'
$ws.Cells.Item(107, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(111, 9).Value = 'sv'
$ws.Cells.Item(111, 10).Value = 'Statement-opinion'
$ws.Cells.Item(114, 9).Value = 'sv'
$ws.Cells.Item(114, 10).Value = 'Statement-opinion'
$ws.Cells.Item(119, 9).Value = 'sd'
$ws.Cells.Item(119, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(128, 9).Value = 'ba'
$ws.Cells.Item(128, 10).Value = 'Appreciation'
$ws.Cells.Item(133, 9).Value = 'aa'
$ws.Cells.Item(133, 10).Value = 'Agree/Accept'
$ws.Cells.Item(134, 9).Value = 'aa'
$ws.Cells.Item(134, 10).Value = 'Agree/Accept'
$ws.Cells.Item(146, 9).Value = 'aa'
$ws.Cells.Item(146, 10).Value = 'Agree/Accept'
$ws.Cells.Item(152, 9).Value = 'b'
$ws.Cells.Item(152, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(160, 9).Value = 'aa'
$ws.Cells.Item(160, 10).Value = 'Agree/Accept'
$ws.Cells.Item(161, 9).Value = 'sd'
$ws.Cells.Item(161, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(168, 9).Value = 'sd'
$ws.Cells.Item(168, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(172, 9).Value = 'sd'
$ws.Cells.Item(172, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(190, 9).Value = 'sv'
$ws.Cells.Item(190, 10).Value = 'Statement-opinion'
$ws.Cells.Item(195, 9).Value = 'aa'
$ws.Cells.Item(195, 10).Value = 'Agree/Accept'
$ws.Cells.Item(196, 9).Value = 'aa'
$ws.Cells.Item(196, 10).Value = 'Agree/Accept'
$ws.Cells.Item(201, 9).Value = 'sd'
$ws.Cells.Item(201, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(202, 9).Value = 'sd'
$ws.Cells.Item(202, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(208, 9).Value = 'sd'
$ws.Cells.Item(208, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(215, 9).Value = 'b'
$ws.Cells.Item(215, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(220, 9).Value = 'aa'
$ws.Cells.Item(220, 10).Value = 'Agree/Accept'
$ws.Cells.Item(221, 9).Value = 'sd'
$ws.Cells.Item(221, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(222, 9).Value = 'aa'
$ws.Cells.Item(222, 10).Value = 'Agree/Accept'
$ws.Cells.Item(223, 9).Value = 'aa'
$ws.Cells.Item(223, 10).Value = 'Agree/Accept'
$ws.Cells.Item(230, 9).Value = 'sd'
$ws.Cells.Item(230, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(235, 9).Value = 'sd'
$ws.Cells.Item(235, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(240, 9).Value = 'sv'
$ws.Cells.Item(240, 10).Value = 'Statement-opinion'
$ws.Cells.Item(255, 9).Value = 'aa'
$ws.Cells.Item(255, 10).Value = 'Agree/Accept'
$ws.Cells.Item(281, 9).Value = 'sd'
$ws.Cells.Item(281, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(297, 9).Value = 'sd'
$ws.Cells.Item(297, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(299, 9).Value = 'sd'
$ws.Cells.Item(299, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(307, 9).Value = 'sv'
$ws.Cells.Item(307, 10).Value = 'Statement-opinion'
$ws.Cells.Item(314, 9).Value = 'ba'
$ws.Cells.Item(314, 10).Value = 'Appreciation'
$ws.Cells.Item(316, 9).Value = 'aa'
$ws.Cells.Item(316, 10).Value = 'Agree/Accept'
$ws.Cells.Item(322, 9).Value = 'aa'
$ws.Cells.Item(322, 10).Value = 'Agree/Accept'
$ws.Cells.Item(325, 9).Value = 'aa'
$ws.Cells.Item(325, 10).Value = 'Agree/Accept'
$ws.Cells.Item(331, 9).Value = 'sd'
$ws.Cells.Item(331, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(337, 9).Value = 'sv'
$ws.Cells.Item(337, 10).Value = 'Statement-opinion'
$ws.Cells.Item(366, 9).Value = 'sv'
$ws.Cells.Item(366, 10).Value = 'Statement-opinion'
$ws.Cells.Item(367, 9).Value = 'sd'
$ws.Cells.Item(367, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(368, 9).Value = '%'
$ws.Cells.Item(368, 10).Value = 'Uninterpretable'
$ws.Cells.Item(376, 9).Value = 'b'
$ws.Cells.Item(376, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(380, 9).Value = 'sv'
$ws.Cells.Item(380, 10).Value = 'Statement-opinion'
$ws.Cells.Item(381, 9).Value = 'sd'
$ws.Cells.Item(381, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(391, 9).Value = '%'
$ws.Cells.Item(391, 10).Value = 'Uninterpretable'
$ws.Cells.Item(398, 9).Value = 'sd'
$ws.Cells.Item(398, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(402, 9).Value = 'sd'
$ws.Cells.Item(402, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(408, 9).Value = 'sd'
$ws.Cells.Item(408, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(418, 9).Value = 'sd'
$ws.Cells.Item(418, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(427, 9).Value = 'sv'
$ws.Cells.Item(427, 10).Value = 'Statement-opinion'
$ws.Cells.Item(441, 9).Value = '%'
$ws.Cells.Item(441, 10).Value = 'Uninterpretable'
$ws.Cells.Item(444, 9).Value = 'sd'
$ws.Cells.Item(444, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(466, 9).Value = 'sd'
$ws.Cells.Item(466, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(472, 9).Value = 'sd'
$ws.Cells.Item(472, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(473, 9).Value = 'ba'
$ws.Cells.Item(473, 10).Value = 'Appreciation'
$ws.Cells.Item(486, 9).Value = 'aa'
$ws.Cells.Item(486, 10).Value = 'Agree/Accept'
$ws.Cells.Item(490, 9).Value = 'qy'
$ws.Cells.Item(490, 10).Value = 'Yes-No-Question'
$ws.Cells.Item(492, 9).Value = 'sd'
$ws.Cells.Item(492, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(495, 9).Value = 'sv'
$ws.Cells.Item(495, 10).Value = 'Statement-opinion'
$ws.Cells.Item(503, 9).Value = 'sd'
$ws.Cells.Item(503, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(504, 9).Value = 'sv'
$ws.Cells.Item(504, 10).Value = 'Statement-opinion'
